$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 15.67790914277992
$ws.Range("A3").Value = 19.88800271446007
$ws.Range("A4").Value = 18.89080234373404
$ws.Range("A5").Value = 20.28454590844126
$ws.Range("A6").Value = 10.55267256212903
$ws.Range("A7").Value = 20.63310720507212
$ws.Range("A8").Value = 16.58696707729638
$ws.Range("A9").Value = 15.38926229123268
$ws.Range("A10").Value = 21.34770284898741
$ws.Range("A11").Value = 22.6759862887024
$ws.Range("A12").Value = 13.13312222896269
$ws.Range("A13").Value = 18.45447784467422
$ws.Range("A14").Value = 16.64593431388198
$ws.Range("A15").Value = 8.17536073061018
$ws.Range("A16").Value = 8.947183266864215

# Rows 17-86 share the same value
$ws.Range("A17:A86").Value = 4.398575093383272

